# Weekly sprint sheet rollover: advance the "Week Of" header to the next
# week and shift each action item's Last/This/Next-week notes forward one
# slot, filling in the new "Next Week" plans.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New week header (row 3, under "Project: 5" in row 2)
$ws.Range("A3").Value = "Week Of:  April 29 - May 3, 2019"

# Row 9 - "Game Logic" action item (Juan Guiterrez)
$ws.Range("D9").Value  = "Start deciding how game will be played out"
$ws.Range("E9").Value  = "Worked around with map border and images"
$ws.Range("F9").Value  = "Add obstacles and power ups"

# Row 10 - "Configuring game scenes..." action item (Jordan Chen)
$ws.Range("D10").Value = "Decide how Tank will function in game"
$ws.Range("E10").Value = "Worked around with game bonuses"
$ws.Range("F10").Value = "debug"

# Row 11 - "Networking features..." action item (Edgar Camacho)
$ws.Range("D11").Value = "Start Server for Game"
$ws.Range("E11").Value = "Write Game Logic Code"
$ws.Range("F11").Value = "Game appearance"

# Row 12 - "GUI elements..." action item (Bryan Nguyen)
$ws.Range("D12").Value = "Create ideas for how the GUI and user interface will look like"
$ws.Range("E12").Value = "Worked with player and bullet objects"
$ws.Range("F12").Value = "Lobby and Start-up appearance"

# Reflect the author's last selection/cursor position on the sheet.
$null = $ws.Range("F10").Select()
